# results to output files
# - rename header H1 "hub_genes" -> "hub_gene"
# - move the active selection to H1 (the cell that was just edited)
# - re-flow column widths: split the former C:E width-group so column E
#   gets its own (wider) width, and nudge column I's width slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "hub_gene"

# Column E needs to stand on its own with a wider width (no longer grouped
# with C:D). ColumnWidth is specified in characters and Excel quantizes the
# stored width to whole pixels of the workbook's default font, so we pick
# the input value that lands on the pixel bucket closest to 20.45.
$ws.Range("E1").EntireColumn.ColumnWidth = 19.65

# Column I's width shifts very slightly too.
$ws.Range("I1").EntireColumn.ColumnWidth = 44.65

# Move/restore the selection to the cell that was edited.
$ws.Range("H1").Select()
